$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range('D2').Value = '60.127.38'
$ws.Range('E2').Value = '  +1.52%  '

$ws.Range('D3').Value = '2.588.62'
$ws.Range('E3').Value = '  +0.22%  '

Set-TextCell $ws 'D5' '564.88'
$ws.Range('E5').Value = '  +0.25%  '

Set-TextCell $ws 'D6' '141.51'
$ws.Range('E6').Value = '  -0.84%  '

Set-TextCell $ws 'D7' '0.997'
$ws.Range('E7').Value = '  -0.14%  '

Set-TextCell $ws 'D8' '0.596'
$ws.Range('E8').Value = '  -0.50%  '

$ws.Range('D9').Value = '2.606.72'
$ws.Range('E9').Value = '  +0.61%  '

Set-TextCell $ws 'D10' '6.56'
$ws.Range('E10').Value = '  -1.26%  '

$ws.Range('E11').Value = '  +0.81%  '

Set-TextCell $ws 'D12' '0.368'
$ws.Range('E12').Value = '  +6.08%  '

Set-TextCell $ws 'D13' '0.150'
$ws.Range('E13').Value = '  -6.27%  '

$ws.Range('D14').Value = '3.048.49'
$ws.Range('E14').Value = '  +0.36%  '

$ws.Range('D15').Value = '60.136.87'
$ws.Range('E15').Value = '  +1.60%  '

Set-TextCell $ws 'D16' '23.22'
$ws.Range('E16').Value = '  +1.61%  '

$ws.Range('E17').Value = '  +1.72%  '

$ws.Range('D18').Value = '2.597.80'
$ws.Range('E18').Value = '  +0.49%  '

Set-TextCell $ws 'D19' '11.27'
$ws.Range('E19').Value = '  +8.56%  '

Set-TextCell $ws 'D20' '4.64'
$ws.Range('E20').Value = '  +1.73%  '

Set-TextCell $ws 'D21' '345.16'
$ws.Range('E21').Value = '  +2.31%  '

Set-TextCell $ws 'D22' '6.92'
$ws.Range('E22').Value = '  +7.84%  '

Set-TextCell $ws 'D23' '1.00'
$ws.Range('E23').Value = '  -0.09%  '

Set-TextCell $ws 'D24' '0.536'
$ws.Range('E24').Value = '  +16.73%  '

Set-TextCell $ws 'D25' '62.90'
$ws.Range('E25').Value = '  -2.19%  '

Set-TextCell $ws 'D26' '0.996'
$ws.Range('E26').Value = '  -0.24%  '

$ws.Range('E27').Value = '  -2.17%  '

Set-TextCell $ws 'D28' '7.61'
$ws.Range('E28').Value = '  +4.02%  '

$ws.Range('D29').Value = '0.0₃0779'
$ws.Range('E29').Value = '  +0.43%  '

$ws.Range('E30').Value = '  +6.55%  '

$ws.Range('E31').Value = '  -0.08%  '

Set-TextCell $ws 'D32' '6.30'
$ws.Range('E32').Value = '  +3.08%  '

Set-TextCell $ws 'D33' '161.14'
$ws.Range('E33').Value = '  -0.15%  '

Set-TextCell $ws 'D34' '19.39'
$ws.Range('E34').Value = '  +2.44%  '

Set-TextCell $ws 'D35' '4.21'
$ws.Range('E35').Value = '  +4.77%  '

$ws.Range('E36').Value = '  +8.47%  '

$ws.Range('E37').Value = '  +3.64%  '

Set-TextCell $ws 'D38' '1.59'
$ws.Range('E38').Value = '  +7.15%  '

Set-TextCell $ws 'D39' '37.68'
$ws.Range('E39').Value = '  +0.51%  '

$ws.Range('E40').Value = '  -2.36%  '

Set-TextCell $ws 'D41' '3.80'
$ws.Range('E41').Value = '  +3.93%  '

Set-TextCell $ws 'D42' '292.71'
$ws.Range('E42').Value = '  -0.84%  '

Set-TextCell $ws 'D43' '138.03'
$ws.Range('E43').Value = '  +4.49%  '

$ws.Range('E44').Value = '  -0.27%  '

Set-TextCell $ws 'D45' '0.0977'
$ws.Range('E45').Value = '  +0.32%  '

Set-TextCell $ws 'D46' '0.603'
$ws.Range('E46').Value = '  +1.13%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D47' '19.45'
$ws.Range('E47').Value = '  +1.93%  '

$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D48' '0.0542'
$ws.Range('E48').Value = '  +1.02%  '

$ws.Range('E49').Value = '  +2.31%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws 'D50' '4.88'
$ws.Range('E50').Value = '  +8.14%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell $ws 'D51' '10.65'
$ws.Range('E51').Value = '  +0.05%  '
